$wb = $excel.ActiveWorkbook

# --- Update Sheet2's selection (was C14, now A7:C11 w/ active cell C11) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A7:C11").Select()

# --- Add Sheet3 after Sheet2, becomes the new active sheet/tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# --- Populate Sheet3 with the new gem-taking summary table ---
$ws3.Range("A1").Value = "take"
$ws3.Range("B1").Value = 15
$ws3.Range("C1").Value = "5 take 1, 5 take 2, 5 discard"

$ws3.Range("A2").Value = "Reserve"
$ws3.Range("B2").Value = 12
$ws3.Range("C2").Value = "any of 12"

$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = "top of deck"

$ws3.Range("A4").Value = "buy"
$ws3.Range("B4").Value = 12
$ws3.Range("C4").Value = "any of 12"

$ws3.Range("B5").Value = 3
$ws3.Range("C5").Value = "reserves"

$ws3.Range("B6").Value = 15
$ws3.Range("C6").Value = "with gold"

$ws3.Range("B7").Formula = "=SUM(B1:B6)"

$ws3.Range("B8").Select()
